$d = $word.ActiveDocument

# Locate only the trailing "):" run (leave the "Revisor (" / "es" runs and
# the gramStart/gramEnd proofErr markers between them untouched) and extend
# its text - this way the new text inherits that run's character
# formatting (sz/szCs = 24) instead of picking up default formatting.
$rng = $d.Content
$found = $rng.Find.Execute("Revisor (es):", $false, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)

if ($found) {
    $tail = $d.Range($rng.End - 2, $rng.End)   # just "):"
    $tail.Text = "): Rafael Souza"

    # Re-apply (identical) direct character formatting to just the appended
    # name so Word splits it back out into its own run - matching the way
    # Word keeps explicitly-applied run formatting in a run of its own -
    # while leaving the original "):" run/properties untouched.
    $nameRng = $d.Range($tail.End - 13, $tail.End)
    $nameRng.Bold = 1
    $nameRng.Bold = 0
}
